# This edit moves the whole sequence-diagram lifeline layout on slide 2 to the
# right by 925158 EMU (~72.85pt), matching the upstream commit that reworked the
# doc with the builder pattern / singleton-pattern diagrams / async-work write-up
# (the participants & messages on this slide were nudged right to make room).
#
# Every top-level shape on the slide gets the same new "Left"; one connector
# ("Straight Connector 4") is glued at its start point (stCxn) to "Rectangle 3",
# so besides Left it also needs its Width/Height updated to keep tracking the
# glue point on the rectangle after the move (PowerPoint normally re-routes such
# connectors automatically when you drag the shape they are glued to).
#
# NOTE: the PowerPoint object model keeps Left/Top/Width/Height in points using
# single-precision (float32) storage internally. A naive EMU/12700 literal can
# therefore round-trip to one EMU less than intended once saved back to OOXML, so
# the point literals below are nudged to the nearest float32 that reproduces the
# exact target EMU on save.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$sh = $s.Shapes.Item("Straight Connector 39")
$sh.Left = 558.5157470703125

$sh = $s.Shapes.Item("Rectangle 3")
$sh.Left = 56.472129821777344

$sh = $s.Shapes.Item("Straight Connector 4")
$sh.Left = 138.07057189941406
$sh.Width = 2.7263782024383545
$sh.Height = 397.9656982421875

$sh = $s.Shapes.Item("Rectangle 5")
$sh.Left = 292.94805908203125

$sh = $s.Shapes.Item("Straight Connector 6")
$sh.Left = 353.5057067871094

$sh = $s.Shapes.Item("Group 9")
$sh.Left = 139.17410278320312

$sh = $s.Shapes.Item("Rectangle 38")
$sh.Left = 516.6613159179688

$sh = $s.Shapes.Item("TextBox 43")
$sh.Left = 255.0088348388672

$sh = $s.Shapes.Item("Group 41")
$sh.Left = 359.8599548339844

$sh = $s.Shapes.Item("Elbow Connector 36")
$sh.Left = 138.20860290527344

$sh = $s.Shapes.Item("Straight Connector 74")
$sh.Left = 190.609619140625

$sh = $s.Shapes.Item("Straight Connector 75")
$sh.Left = 140.03395080566406

$sh = $s.Shapes.Item("TextBox 60")
$sh.Left = 268.5805969238281

$sh = $s.Shapes.Item("Group 62")
$sh.Left = 358.86663818359375

$sh = $s.Shapes.Item("TextBox 79")
$sh.Left = 514.9329833984375

$sh = $s.Shapes.Item("TextBox 80")
$sh.Left = 362.1260070800781

$sh = $s.Shapes.Item("TextBox 81")
$sh.Left = 514.5733032226562

$sh = $s.Shapes.Item("TextBox 82")
$sh.Left = 234.0149688720703

$sh = $s.Shapes.Item("TextBox 83")
$sh.Left = 185.63002014160156

